$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the findByTags mock-data JSON (E3): remove stray literal newline inside the
# "rule" value and drop the erroneous leading escaped quote before the "output" JSON ---
$e3 = @'
{
    "resource": "pets",
    "url": "/pets/findByTags",
    "type": "Params",
    "rule": "[{\"tags\":\"spring-grey\"}]",
    "operationId": "findPetsByTags",
    "output": "{\n  \"category\": {\n    \"id\": 200,\n    \"name\": \"Bulldog\"\n  },\n  \"id\": 201,\n  \"name\": \"Butch\",\n  \"photoUrls\": [\n    \"string\"\n  ],\n  \"status\": \"available\",\n  \"tags\": [\n    {\n      \"id\": 201,\n      \"name\": \"<tags>\"\n    }\n  ]\n}",
    "httpStatusCode": "200",
    "method": "GET",
    "availableParams": [
        {
            "key": "tags",
            "value": "<tags>",
            "parameterType": "QUERY_PARAM"
        }
    ]
}
'@
$ws.Range("E3").Value = $e3

# --- Replace the CSV "Csvson" content (H4) with the new i~ prefixed id references ---
$h4 = @'
id,name,category/id:name,photoUrls,status,tags/id:name
i~201,Butch,i~200:Bulldog,string\|,available,i~201:spring-grey\|
'@
$ws.Range("H4").Value = $h4

# --- Fix the addPet mock-data JSON (E5): correct malformed escaping in "input"/"output" ---
$e5 = @'
{
    "resource": "pets",
    "url": "/pets",
    "type": "Response",
    "operationId": "addPet",
    "input": "{\n  \"category\": {\n    \"id\": 100,\n    \"name\": \"string\" \n  },\n  \"id\": 100,\n  \"name\": \"doggie-1\",\n  \"photoUrls\": [\n    \"string\" \n  ],\n  \"status\": \"available\",\n  \"tags\": [\n    {\n      \"id\": 0,\n      \"name\": \"string\" \n    }\n  ]\n}",
    "output": "{\n  \"category\": {\n    \"id\": 100,\n    \"name\": \"string\" \n  },\n  \"id\": 100,\n  \"name\": \"doggie\",\n  \"photoUrls\": [\n    \"string\" \n  ],\n  \"status\": \"available\",\n  \"tags\": [\n    {\n      \"id\": 0,\n      \"name\": \"string\" \n    }\n  ]\n}",
    "httpStatusCode": "201",
    "method": "POST"
}
'@
$ws.Range("E5").Value = $e5

# --- Update the sheet view: scroll position and active selection ---
$ws.Range("F5").Select()

# --- Row 5 height changed from 280.5 to 267.75 ---
$ws.Rows.Item(5).RowHeight = 267.75
